# Registro etapa 2 Demo - add daily trade registrations for weeks 2 and 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 2 (rows 28-32): add entries for 29-Aug-2025 and 02-Sep-2025
$ws.Range("A31").Value = 45898
$ws.Range("C31").Value = 1
$ws.Range("E31").Value = 0

$ws.Range("A32").Value = 45902
$ws.Range("C32").Value = 1
$ws.Range("E32").Value = 0

# Week 3 (rows 36-40): add entries for 30-Sep-2025, 01-Oct-2025, 07-Oct-2025, 08-Oct-2025
$ws.Range("A36").Value = 45930
$ws.Range("C36").Value = 1
$ws.Range("E36").Value = 0

$ws.Range("A37").Value = 45931
$ws.Range("C37").Value = 1
$ws.Range("E37").Value = 0

$ws.Range("A38").Value = 45937
$ws.Range("C38").Value = 1
$ws.Range("E38").Value = 0

$ws.Range("A39").Value = 45938
$ws.Range("C39").Value = 1
$ws.Range("E39").Value = 0

# Move the active selection to the latest entry row (E39:F39)
$ws.Activate()
$ws.Range("E39:F39").Select()
